$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "StartEffects" column (I) was always empty - it's removed entirely,
# shifting bIsRanged/DropExpValue/ProjectileDataClass one column to the left.
$ws.Columns("I").Delete()

# Re-apply the custom widths for the columns that now sit at I..L
# (9 -> raw 14, 10 -> raw 16.125 [closest attainable 16.142857], 11 -> raw 23.5 [closest attainable], 12 -> raw 9)
$ws.Columns("I").ColumnWidth = 13.2857142857143
$ws.Columns("J").ColumnWidth = 15.4107142857143
$ws.Columns("K").ColumnWidth = 22.7857142857143
$ws.Columns("L").ColumnWidth = 8.28571428571429

# Header rename: StartAbilities -> EnemyAbilities (now the single ability data-asset column)
$ws.Range("H1").Value = "EnemyAbilities"

# Normal monsters (rows 2-6): ability blueprint list replaced by a reference
# to the shared Normal enemy abilities DataAsset.
$ws.Range("H2:H6").Value = "/Game/Data/Enemy/DA/Normal/DA_NormalEnemyAbilitiesData.DA_NormalEnemyAbilitiesData"

# Boss monsters (rows 7-8): ability blueprint list replaced by a reference
# to the shared Boss enemy abilities DataAsset.
$ws.Range("H7:H8").Value = "/Game/Data/Enemy/DA/Boss/DA_BossEnemyAbilitiesData.DA_BossEnemyAbilitiesData"

# Normal_Lich gains a dedicated ProjectileDataClass value (was IceSpear placeholder, now Lich's own attack BP).
$ws.Range("K3").Value = "/Game/Personal/LIM_H_S/BP_LichAttack.BP_LichAttack_C"

# Boss_Sevarog mesh path swapped from the Chronos skin to the base Sevarog mesh.
$ws.Range("B7").Value = "/Game/Aseets/ParagonSevarog/Characters/Heroes/Sevarog/Meshes/Sevarog.Sevarog"

# Move the active selection, matching the author's final cursor position.
$ws.Range("H17").Select()
